$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rewrite "Curso Basico" step 4 text (row 19) - drop the login-validation clause,
#    then delete the three login sub-steps above it (rows 16-18) and the two
#    alternate-course sub-steps below the Curso Alternativo row (rows 22-23), plus
#    blank out the Curso Alternativo text itself (row 21) - mirrors removing the
#    basic "log in to the system" course from the use case.
$ws.Range("B19").Value = "El sistema muestra el formulario Administración de supervisores la cual presenta:`n*Listado de telemarketers asignados al mismo grupo que el supervisor, con botón `"Ver`" en cada uno de los registros.`n*Listado de oportunidades asignadas al mismo grupo que el supervisor, con botón `"Ver`" en cada uno de los registros.`n*Listado de campañas asignadas al mismo grupo que el supervisor, con botón `"Ver`" en cada uno de los registros.`n*Listado de ventas asignadas al mismo grupo que el supervisor, con botón `"Ver`" en cada uno de los registros.`n"
$ws.Range("B21").ClearContents()

# 2) Update 'paso 3.x' -> 'paso 2.x' references (rows 11-14)
$ws.Range("B11").Value = "Ver Telemarketer, paso 2.a"
$ws.Range("B12").Value = "Ver Oportunidad, paso 2.b"
$ws.Range("B13").Value = "Ver Campaña, paso 2.c"
$ws.Range("B14").Value = "Ver Venta, paso 2.d"

# 3) Update the Pre-condicion text (row 9)
$ws.Range("B9").Value = "El usuario debe pertenecer al tipo de Usuario `"Supervisor`"."

# Delete rows bottom-up so the earlier row numbers stay valid while deleting.
$ws.Rows("22:23").Delete()
$ws.Rows("16:18").Delete()

# 4) Renumber the remaining "Curso Basico" steps (were 4 and 5, now 1 and 2)
$ws.Range("A16").Value = 1
$ws.Range("A17").Value = 2

# 5) Restore the active selection to B9 (matches the saved view state in the diff)
$ws.Range("B9").Select()

